$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1686746987951807
$ws.Range("C2").Value = 0.5963855421686747
$ws.Range("J2").Value = 0.01506024096385542
$ws.Range("P2").Value = 0.1385542168674699
$ws.Range("S2").Value = 0.08132530120481928
$ws.Range("B3").Value = 0.01005025125628141
$ws.Range("C3").Value = 0.02010050251256281
$ws.Range("J3").Value = 0.02512562814070352
$ws.Range("P3").Value = 0.6834170854271356
$ws.Range("S3").Value = 0.2613065326633166
$ws.Range("J4").Value = 0.03773584905660377
$ws.Range("P4").Value = 0.6415094339622641
$ws.Range("S4").Value = 0.3207547169811321
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.1238938053097345
$ws.Range("D6").Value = 0.02654867256637168
$ws.Range("F6").Value = 0.07079646017699115
$ws.Range("J6").Value = 0.247787610619469
$ws.Range("O6").Value = 0.03097345132743363
$ws.Range("Q6").Value = 0.1725663716814159
$ws.Range("R6").Value = 0.02654867256637168
$ws.Range("S6").Value = 0.3008849557522124
$ws.Range("B7").Value = 0.1318681318681319
$ws.Range("D7").Value = 0.01648351648351648
$ws.Range("F7").Value = 0.06043956043956044
$ws.Range("J7").Value = 0.1373626373626374
$ws.Range("O7").Value = 0.03846153846153846
$ws.Range("Q7").Value = 0.1098901098901099
$ws.Range("R7").Value = 0.06043956043956044
$ws.Range("S7").Value = 0.445054945054945
$ws.Range("B8").Value = 0.1271820448877805
$ws.Range("D8").Value = 0.02743142144638404
$ws.Range("F8").Value = 0.06234413965087282
$ws.Range("J8").Value = 0.1047381546134663
$ws.Range("O8").Value = 0.01745635910224439
$ws.Range("Q8").Value = 0.1820448877805486
$ws.Range("R8").Value = 0.08728179551122195
$ws.Range("S8").Value = 0.3915211970074813
$ws.Range("B9").Value = 0.09905660377358491
$ws.Range("D9").Value = 0.0330188679245283
$ws.Range("F9").Value = 0.07547169811320754
$ws.Range("J9").Value = 0.1132075471698113
$ws.Range("Q9").Value = 0.1886792452830189
$ws.Range("R9").Value = 0.08962264150943396
$ws.Range("S9").Value = 0.4009433962264151
$ws.Range("B10").Value = 0.1258680555555556
$ws.Range("D10").Value = 0.0234375
$ws.Range("E10").Value = 0.003472222222222222
$ws.Range("F10").Value = 0.07118055555555555
$ws.Range("J10").Value = 0.1284722222222222
$ws.Range("O10").Value = 0.01909722222222222
$ws.Range("Q10").Value = 0.1961805555555556
$ws.Range("R10").Value = 0.06944444444444445
$ws.Range("S10").Value = 0.3628472222222222
$ws.Range("G11").Value = 0.1928571428571429
$ws.Range("J11").Value = 0.075
$ws.Range("K11").Value = 0.2392857142857143
$ws.Range("L11").Value = 0.475
$ws.Range("S11").Value = 0.01785714285714286
$ws.Range("G12").Value = 0.7857142857142857
$ws.Range("J12").Value = 0.1357142857142857
$ws.Range("L12").Value = 0.03571428571428571
$ws.Range("S12").Value = 0.04285714285714286
$ws.Range("G13").Value = 0.6176470588235294
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.1470588235294118
$ws.Range("F15").Value = 0.02834008097165992
$ws.Range("H15").Value = 0.1417004048582996
$ws.Range("I15").Value = 0.1012145748987854
$ws.Range("J15").Value = 0.3441295546558704
$ws.Range("K15").Value = 0.05668016194331984
$ws.Range("M15").Value = 0.008097165991902834
$ws.Range("O15").Value = 0.08502024291497975
$ws.Range("S15").Value = 0.2348178137651822
$ws.Range("F16").Value = 0.02369668246445497
$ws.Range("H16").Value = 0.1611374407582938
$ws.Range("I16").Value = 0.08530805687203792
$ws.Range("J16").Value = 0.3412322274881517
$ws.Range("K16").Value = 0.0995260663507109
$ws.Range("M16").Value = 0.01895734597156398
$ws.Range("N16").Value = 0.004739336492890996
$ws.Range("O16").Value = 0.1232227488151659
$ws.Range("S16").Value = 0.1421800947867299
$ws.Range("F17").Value = 0.02290076335877863
$ws.Range("H17").Value = 0.1679389312977099
$ws.Range("I17").Value = 0.0916030534351145
$ws.Range("J17").Value = 0.4173027989821883
$ws.Range("K17").Value = 0.08905852417302799
$ws.Range("M17").Value = 0.01272264631043257
$ws.Range("O17").Value = 0.08651399491094147
$ws.Range("S17").Value = 0.1119592875318066
$ws.Range("F18").Value = 0.0392156862745098
$ws.Range("H18").Value = 0.1568627450980392
$ws.Range("I18").Value = 0.1372549019607843
$ws.Range("J18").Value = 0.3856209150326798
$ws.Range("K18").Value = 0.0915032679738562
$ws.Range("M18").Value = 0.0130718954248366
$ws.Range("O18").Value = 0.0457516339869281
$ws.Range("S18").Value = 0.130718954248366
$ws.Range("F19").Value = 0.02238805970149254
$ws.Range("H19").Value = 0.1956882255389718
$ws.Range("I19").Value = 0.09286898839137644
$ws.Range("J19").Value = 0.3665008291873963
$ws.Range("K19").Value = 0.1077943615257048
$ws.Range("M19").Value = 0.01741293532338309
$ws.Range("O19").Value = 0.07048092868988391
$ws.Range("S19").Value = 0.1268656716417911
